$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.180599999999991
$ws.Range("D14").Value = -8.2433
$ws.Range("D16").Value = -8.150599999999995
$ws.Range("D21").Value = -7.864100000000002
$ws.Range("D23").Value = -7.383699999999992
$ws.Range("D25").Value = -8.207000000000001
